$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.411.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = "'1.566.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.41%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D6").Value = "'284.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.37%  '
$ws.Range("E7").Value = '  -2.13%  '
$ws.Range("D8").Value = "'48.56"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.83%  '
$ws.Range("D9").Value = "'0.3328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.01%  '
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("D11").Value = "'0.07395"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.20%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").Value = "'5.948"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("D15").Value = "'6.897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.88%  '
$ws.Range("D16").Value = "'1.564.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("E17").Value = '  -1.70%  '
$ws.Range("D18").Value = "'88.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("D19").Value = "'0.06702"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = "'6.347"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.75%  '
$ws.Range("D22").Value = "'16.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = "'11.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = "'22.399.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = "'2.384"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("D26").Value = "'2.536"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.90%  '
$ws.Range("D27").Value = "'149.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").Value = "'19.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.67%  '
$ws.Range("D29").Value = "'4.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("D30").Value = "'123.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.51%  '
$ws.Range("D31").Value = "'1.742.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.31%  '
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").Value = "'6.100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = "'1.990"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("D35").Value = "'9.792"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").Value = "'0.08273"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -3.02%  '
$ws.Range("D38").Value = "'0.2234"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("D39").Value = "'0.06421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").Value = "'5.362"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.21%  '
$ws.Range("E41").Value = '  -6.91%  '
$ws.Range("D42").Value = "'0.6233"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = "'11.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").Value = "'1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = "'13.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = "'0.6006"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("E47").Value = '  -1.76%  '
$ws.Range("D48").Value = "'2.030"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").Value = "'123.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.41%  '
$ws.Range("D50").Value = "'1.212"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("D51").Value = "'0.07199"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.79%  '
